$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 85: Tyler Morgan-Wall
$ws.Range("A85").Value = "Tyler Morgan-Wall"
$ws.Range("C85").Value = "R"
$ws.Range("E85").Value = "tylermorganwall"
$ws.Range("F85").Value = "https://www.tylermw.com/"
$ws.Hyperlinks.Add($ws.Range("F85"), "https://www.tylermw.com/") | Out-Null

# Row 86: Gavin Simpson
$ws.Range("A86").Value = "Gavin Simpson"
$ws.Range("C86").Value = "R"
$ws.Range("E86").Value = "ucfagls"
$ws.Range("F86").Value = "https://www.fromthebottomoftheheap.net"
$ws.Hyperlinks.Add($ws.Range("F86"), "https://www.fromthebottomoftheheap.net") | Out-Null

$ws.Range("A85:A86").Style = $ws.Range("A84").Style
$ws.Range("F85:F86").Style = $ws.Range("F84").Style

$ws.Range("E89").Select()
